$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53: transaction finalized (now DONE), fill in finalized date & fee ---
$ws.Cells.Item(53,8).Value = "DONE"
$ws.Cells.Item(53,9).Value = 42866.194872685184
$ws.Cells.Item(53,10).Value = "0.26918413 XRP (0.15%)"

# --- Row 54: new transaction row, copy number formats from row 53 first ---
$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("H53").Copy()
$ws.Range("H54").PasteSpecial(-4122)
$ws.Range("I53").Copy()
$ws.Range("I54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# D54 needs to hold a text string that looks like a number, so force text
# format before assigning the value (otherwise it is parsed as a number),
# then restore the wrapped/general style used by the rest of the column.
$ws.Cells.Item(54,4).NumberFormat = "@"
$ws.Cells.Item(54,4).Value = " 0.17838400`n"
$ws.Range("D53").Copy()
$ws.Range("D54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(54,1).Value = 42866.28197916667
$ws.Cells.Item(54,2).Value = "            Sell"
$ws.Cells.Item(54,3).Value = "        XRP"
$ws.Cells.Item(54,5).Value = "         0.185USDT"
$ws.Cells.Item(54,6).Value = "         180 XRP"
$ws.Cells.Item(54,7).Value = " XRP/USDT0000005"
$ws.Cells.Item(54,8).Value = "IN PROGRESS"
$ws.Cells.Item(54,11).Value = "     "

# Keep the same row height as the rest of the table (Excel would otherwise
# auto-fit the wrapped cell containing an embedded line break).
$ws.Rows(54).RowHeight = 14.25

$ws.Range("C63").Select() | Out-Null
